$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 160 (shifts the existing rows 160-210 down to 161-211)
$ws.Rows.Item(160).Insert()

# Populate the new row 160 with the latest weekly entry
$ws.Cells.Item(160, 1).Value = 11
$ws.Cells.Item(160, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(160, 3).Value = "Bíobío"
$ws.Cells.Item(160, 4).Value = 44588
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 100112009
$ws.Cells.Item(160, 7).Value = "Acelga"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 200
$ws.Cells.Item(160, 11).Value = 600
$ws.Cells.Item(160, 12).Value = 650
$ws.Cells.Item(160, 13).Value = 625
$ws.Cells.Item(160, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(160, 15).Value = "Región de Ñuble"
$ws.Cells.Item(160, 16).Value = 625
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# Keep D160 formatted the same way as the date column below it
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(161, 4).NumberFormat
